$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C15").Value = "'94"
$ws.Range("D15").Value = "'230152.38"
$ws.Range("C16").Value = "'408"
$ws.Range("D16").Value = "'1229526.14"
$ws.Range("C17").Value = "'145"
$ws.Range("D17").Value = "'621119.72"
$ws.Range("C18").Value = "'43"
$ws.Range("D18").Value = "'206045.00"
$ws.Range("C19").Value = "'15"
$ws.Range("D19").Value = "'99716.00"
$ws.Range("C20").Value = "'15"
$ws.Range("D20").Value = "'32621.00"
$ws.Range("C23").Value = "'111"
$ws.Range("D23").Value = "'425337.00"
$ws.Range("C33").Value = "'91"
$ws.Range("D33").Value = "'239826.00"
$ws.Range("C34").Value = "'499"
$ws.Range("D34").Value = "'1499825.82"
$ws.Range("C35").Value = "'199"
$ws.Range("D35").Value = "'944847.11"
$ws.Range("C39").Value = "'32"
$ws.Range("D39").Value = "'78330.00"
$ws.Range("C40").Value = "'151"
$ws.Range("D40").Value = "'384328.00"
$ws.Range("C41").Value = "'76"
$ws.Range("D41").Value = "'262900.00"
$ws.Range("C42").Value = "'19"
$ws.Range("D42").Value = "'75995.14"
$ws.Range("C44").Value = "'41"
$ws.Range("D44").Value = "'94683.00"
$ws.Range("C50").Value = "'90"
$ws.Range("D50").Value = "'249228.17"
$ws.Range("C51").Value = "'524"
$ws.Range("D51").Value = "'1672527.52"
$ws.Range("C52").Value = "'242"
$ws.Range("D52").Value = "'979735.76"
$ws.Range("C53").Value = "'80"
$ws.Range("D53").Value = "'452878.23"
$ws.Range("C56").Value = "'667"
$ws.Range("D56").Value = "'1668796.41"
$ws.Range("C57").Value = "'3280"
$ws.Range("D57").Value = "'9822525.10"
$ws.Range("C58").Value = "'1693"
$ws.Range("D58").Value = "'6733611.92"
$ws.Range("C59").Value = "'577"
$ws.Range("D59").Value = "'2712640.96"
$ws.Range("C60").Value = "'118"
$ws.Range("D60").Value = "'804041.00"
$ws.Range("C62").Value = "'272"
$ws.Range("D62").Value = "'638263.00"
$ws.Range("C79").Value = "'219"
$ws.Range("D79").Value = "'555826.09"
$ws.Range("C80").Value = "'843"
$ws.Range("D80").Value = "'2576291.11"
$ws.Range("C81").Value = "'316"
$ws.Range("D81").Value = "'1220440.79"
$ws.Range("C82").Value = "'106"
$ws.Range("D82").Value = "'492484.52"
$ws.Range("C91").Value = "'92"
$ws.Range("D91").Value = "'221878.00"
$ws.Range("C92").Value = "'389"
$ws.Range("D92").Value = "'1112260.67"
$ws.Range("C96").Value = "'12"
$ws.Range("D96").Value = "'24000.00"
